$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Version: 1.0 -> 1.2.5
$ws.Range("D2").Value = "1.2.5"

# 2. Precondition text: fix "usuario" -> "usuário" and add trailing period
$precondition = "O usuário devidamente autenticado e na tela inicial de cancelar diárias."
$ws.Range("B8").Value = $precondition
$ws.Range("B17").Value = $precondition
$ws.Range("B25").Value = $precondition

# 3. MSG102 text: add trailing period
$msg102 = "SYSTEM Exibe a mensagem (MSG102 - Confirmar cancelamento)."
$ws.Range("D10").Value = $msg102
$ws.Range("D19").Value = $msg102

# 4. MSG217 text: remove stray tab character before closing parenthesis
$ws.Range("D20").Value = "SYSTEM Identifica que o usuário não informou uma justificativa para o cancelamento. Não efetiva o cancelamento e exibe mensagem de erro (MSG217 - Necessário informar uma justificativa para o cancelamento de solicitações) para o usuário."

# 5. MSG205 text: fix "Solcitação" -> "Solicitação"
$ws.Range("D27").Value = "SYSTEM Identifica que a solicitação de diária está em situação diferente de 'SOLICITADA PARA EMPENHO' ou 'SOLICITADA PARA PRESTAÇÃO DE CONTAS'.  Impede o cancelamento e exibe mensagem de erro (MSG205 - Solicitação de diária não pode ser cancelada) para o usuário."
